# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.860.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.665.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.37%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.401"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.39%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.144.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.703.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.664.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "355.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("E24").Value = "  +9.80%  "

$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.50%  "

$ws.Range("E27").Value = "  +2.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "571.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.02%  "

$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.47%  "

$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("E38").Value = "  +2.70%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0618"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.55%  "

$ws.Range("E46").Value = "  +1.21%  "

$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("E49").Value = "  -1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.818"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.17%  "
